$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-5), reflecting a cyclic shift of the data rows
# (row2 <- old row5, row3 <- old row2, row4 <- old row3, row5 <- old row4)
$data = @(
    @{ Row = 2; D = 44370; I = "Segunda"; J = 100; K = 1000; L = 1200; M = 1080; N = "`$/docena de matas"; P = 180; Q = 6 },
    @{ Row = 3; D = 44623; I = "Primera"; J = 300; K = 1800; L = 2000; M = 1900; N = "`$/paquete";         P = 1900; Q = 1 },
    @{ Row = 4; D = 44377; I = "Segunda"; J = 550; K = 2000; L = 2800; M = 2364; N = "`$/docena de matas"; P = 394; Q = 6 },
    @{ Row = 5; D = 44267; I = "Primera"; J = 120; K = 1500; L = 1800; M = 1650; N = "`$/docena de matas"; P = 275; Q = 6 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value  = $entry.D   # D - Fecha
    $ws.Cells.Item($r, 9).Value  = $entry.I   # I - Calidad
    $ws.Cells.Item($r, 10).Value = $entry.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $entry.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $entry.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $entry.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $entry.N   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $entry.P   # P - Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $entry.Q   # Q - Kg o Unidades
}
